# Refreshed Price (D) and Volume(1h) (E) columns of the cryptos table
# with the latest scrape values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '62.891.40'
$ws.Range("E2").Value = '  +4.94%  '
$ws.Range("D3").Value = '3.115.63'
$ws.Range("E3").Value = '  +3.02%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'585.40"
$ws.Range("E5").Value = '  +3.32%  '
$ws.Range("E6").Value = '  +2.38%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.107.51'
$ws.Range("E8").Value = '  +3.09%  '
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("E10").Value = '  +11.16%  '
$ws.Range("D11").Value = "'5.71"
$ws.Range("E11").Value = '  +7.23%  '
$ws.Range("D12").Value = "'0.468"
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("E13").Value = '  +6.00%  '
$ws.Range("D14").Value = "'35.52"
$ws.Range("E14").Value = '  +3.41%  '
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("D16").Value = '3.634.23'
$ws.Range("E16").Value = '  +3.24%  '
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").Value = '3.112.98'
$ws.Range("E18").Value = '  +2.99%  '
$ws.Range("D19").Value = '62.842.57'
$ws.Range("E19").Value = '  +4.92%  '
$ws.Range("D20").Value = "'465.01"
$ws.Range("E20").Value = '  +6.02%  '
$ws.Range("D21").Value = "'14.10"
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").Value = "'7.53"
$ws.Range("E23").Value = '  +5.33%  '
$ws.Range("D24").Value = "'13.35"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = "'82.17"
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("E28").Value = '  +4.65%  '
$ws.Range("D29").Value = "'8.28"
$ws.Range("E29").Value = '  +5.09%  '
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("D31").Value = "'6.82"
$ws.Range("E31").Value = '  +8.22%  '
$ws.Range("D32").Value = "'26.96"
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("E33").Value = '  +8.23%  '
$ws.Range("D34").Value = '0.0₃0848'
$ws.Range("E34").Value = '  +6.96%  '
$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = '  +11.09%  '
$ws.Range("E36").Value = '  +3.77%  '
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("D38").Value = "'3.25"
$ws.Range("E38").Value = '  +15.80%  '
$ws.Range("D39").Value = "'51.03"
$ws.Range("E39").Value = '  +3.69%  '
$ws.Range("D40").Value = "'432.34"
$ws.Range("E40").Value = '  +6.16%  '
$ws.Range("D41").Value = "'8.79"
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").Value = '2.930.52'
$ws.Range("E43").Value = '  +3.90%  '
$ws.Range("E44").Value = '  +9.02%  '
$ws.Range("E45").Value = '  +3.14%  '
$ws.Range("D46").Value = "'2.18"
$ws.Range("E46").Value = '  +6.54%  '
$ws.Range("D47").Value = "'35.26"
$ws.Range("E47").Value = '  +3.41%  '
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").Value = "'123.17"
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = "'24.75"
$ws.Range("E51").Value = '  +4.44%  '
